# Sync attendance_reports: normalize ordering of "Recorded By" (column G)
# entries that were produced with inconsistent ordering.
#
# Exact, deterministic string substitutions applied to column G cells:
#   "backup@backdoor.com, system, System" -> "backup@backdoor.com, System, system"
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$replacements = @{
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    $value = $cell.Value2

    if ($null -ne $value -and $replacements.ContainsKey($value)) {
        $cell.Value2 = $replacements[$value]
    }
}
